# Update the "想去人数" (F column) counts on the "展览" and "全部类型"
# worksheets to reflect the newly scraped figures.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 3105
    5  = 2728
    7  = 143
    9  = 1445
    13 = 1224
    15 = 368
    19 = 110
    21 = 91
    22 = 2658
    24 = 310
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
